$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109..149 down to
# 110..150. This is a new weekly Cilantro price observation for
# "Terminal La Palmera de La Serena" / Coquimbo.
$ws.Rows("109").Insert()

$ws.Cells.Item(109, 1).Value = 8
$ws.Cells.Item(109, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(109, 3).Value = "Coquimbo"
$ws.Cells.Item(109, 4).Value = 44726
$ws.Cells.Item(109, 5).Value = 4
$ws.Cells.Item(109, 6).Value = 100112040
$ws.Cells.Item(109, 7).Value = "Cilantro"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 3400
$ws.Cells.Item(109, 11).Value = 1500
$ws.Cells.Item(109, 12).Value = 2000
$ws.Cells.Item(109, 13).Value = 1750
$ws.Cells.Item(109, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(109, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(109, 16).Value = 1167
$ws.Cells.Item(109, 17).Value = 1.5
$ws.Cells.Item(109, 18).Value = "Hortaliza"
